$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Cells.Item(2, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "27.539.54"
$cD.Style = $styleD

$cE = $ws.Cells.Item(2, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.42%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(3, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.871.89"
$cD.Style = $styleD

$cE = $ws.Cells.Item(3, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.54%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(4, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.014"
$cD.Style = $styleD

$cE = $ws.Cells.Item(4, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +0.59%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(5, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "312.78"
$cD.Style = $styleD

$cE = $ws.Cells.Item(5, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.09%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(6, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.013"
$cD.Style = $styleD

$cE = $ws.Cells.Item(6, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +0.68%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(7, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.4781"
$cD.Style = $styleD

$cE = $ws.Cells.Item(7, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +0.49%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(8, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.3776"
$cD.Style = $styleD

$cE = $ws.Cells.Item(8, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +3.01%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(9, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.07358"
$cD.Style = $styleD

$cE = $ws.Cells.Item(9, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.25%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(10, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.9377"
$cD.Style = $styleD

$cE = $ws.Cells.Item(10, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.25%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(11, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "20.74"
$cD.Style = $styleD

$cE = $ws.Cells.Item(11, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +5.48%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(12, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.07849"
$cD.Style = $styleD

$cE = $ws.Cells.Item(12, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.01%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(13, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.908.98"
$cD.Style = $styleD

$cE = $ws.Cells.Item(13, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +0.46%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(14, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "5.447"
$cD.Style = $styleD

$cE = $ws.Cells.Item(14, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.57%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(15, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "6.595"
$cD.Style = $styleD

$cE = $ws.Cells.Item(15, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +3.10%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(16, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "90.87"
$cD.Style = $styleD

$cE = $ws.Cells.Item(16, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.44%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(17, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.014"
$cD.Style = $styleD

$cE = $ws.Cells.Item(17, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +0.49%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(18, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.000008923"
$cD.Style = $styleD

$cE = $ws.Cells.Item(18, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +3.26%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(19, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.013"
$cD.Style = $styleD

$cE = $ws.Cells.Item(19, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +0.61%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(20, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "14.98"
$cD.Style = $styleD

$cE = $ws.Cells.Item(20, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.94%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(21, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "27.563.20"
$cD.Style = $styleD

$cE = $ws.Cells.Item(21, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.40%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(22, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "5.143"
$cD.Style = $styleD

$cE = $ws.Cells.Item(22, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.82%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(23, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "10.75"
$cD.Style = $styleD

$cE = $ws.Cells.Item(23, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.08%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(24, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.964"
$cD.Style = $styleD

$cE = $ws.Cells.Item(24, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.39%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(25, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "154.26"
$cD.Style = $styleD

$cE = $ws.Cells.Item(25, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.29%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(26, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "18.55"
$cD.Style = $styleD

$cE = $ws.Cells.Item(26, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.23%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(27, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.024"
$cD.Style = $styleD

$cE = $ws.Cells.Item(27, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.21%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(28, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "115.98"
$cD.Style = $styleD

$cE = $ws.Cells.Item(28, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.61%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(29, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "5.003"
$cD.Style = $styleD

$cE = $ws.Cells.Item(29, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.68%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(30, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.08931"
$cD.Style = $styleD

$cE = $ws.Cells.Item(30, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +0.72%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(31, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "3.342"
$cD.Style = $styleD

$cE = $ws.Cells.Item(31, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.32%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(32, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.219"
$cD.Style = $styleD

$cE = $ws.Cells.Item(32, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +3.98%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(33, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.7564"
$cD.Style = $styleD

$cE = $ws.Cells.Item(33, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.25%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(34, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "4.613"
$cD.Style = $styleD

$cE = $ws.Cells.Item(34, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +3.03%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(35, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.698"
$cD.Style = $styleD

$cE = $ws.Cells.Item(35, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  -1.09%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(36, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.02053"
$cD.Style = $styleD

$cE = $ws.Cells.Item(36, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +5.11%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(37, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.118"
$cD.Style = $styleD

$cE = $ws.Cells.Item(37, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.35%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(38, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.05287"
$cD.Style = $styleD

$cE = $ws.Cells.Item(38, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +0.57%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(39, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "2.999"
$cD.Style = $styleD

$cE = $ws.Cells.Item(39, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.05%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(40, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.5361"
$cD.Style = $styleD

$cE = $ws.Cells.Item(40, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +3.09%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(41, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "7.087"
$cD.Style = $styleD

$cE = $ws.Cells.Item(41, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.00%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(42, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.1527"
$cD.Style = $styleD

$cE = $ws.Cells.Item(42, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.20%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(43, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "8.453"
$cD.Style = $styleD

$cE = $ws.Cells.Item(43, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +3.05%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(44, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "10.70"
$cD.Style = $styleD

$cE = $ws.Cells.Item(44, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.36%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(45, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.4822"
$cD.Style = $styleD

$cE = $ws.Cells.Item(45, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.07%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(46, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.014"
$cD.Style = $styleD

$cE = $ws.Cells.Item(46, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +0.70%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(47, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "1.662"
$cD.Style = $styleD

$cE = $ws.Cells.Item(47, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +3.93%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(48, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "102.82"
$cD.Style = $styleD

$cE = $ws.Cells.Item(48, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.47%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(49, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "67.41"
$cD.Style = $styleD

$cE = $ws.Cells.Item(49, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +2.52%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(50, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.06091"
$cD.Style = $styleD

$cE = $ws.Cells.Item(50, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +1.18%  "
$cE.Style = $styleE

$cD = $ws.Cells.Item(51, 4)
$styleD = $cD.Style
$cD.NumberFormat = "@"
$cD.Value = "0.9263"
$cD.Style = $styleD

$cE = $ws.Cells.Item(51, 5)
$styleE = $cE.Style
$cE.NumberFormat = "@"
$cE.Value = "  +4.71%  "
$cE.Style = $styleE
